# =====================================================================
# Apply the "Updated Files of CF" edit:
#  - Insert two new worksheets: "EmailBrandCampaign" and "GRLocationCampaign"
#  - Populate them with header + one data row each
#  - Update the EmailLocationCampaign sheet's S2 value + selection
#  - Leave Date / Reschedule Date / Sheet1 sheets as-is (content unchanged)
# =====================================================================

$wb = $excel.ActiveWorkbook

$emailLoc = $wb.Worksheets.Item("EmailLocationCampaign")

# ---------------------------------------------------------------------
# 1. Insert the two new worksheets right after EmailLocationCampaign
# ---------------------------------------------------------------------
$brand = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $emailLoc)
$brand.Name = "EmailBrandCampaign"

$grLoc = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $brand)
$grLoc.Name = "GRLocationCampaign"

# ---------------------------------------------------------------------
# 2. EmailBrandCampaign sheet content
# ---------------------------------------------------------------------
# Header row (bold, wrap text)
$brand.Range("A1").Value = "CamType"
$brand.Range("B1").Value = "CamOption"
$brand.Range("C1").Value = "CamLang"
$brand.Range("D1").Value = "CamName"
$brand.Range("E1").Value = "CamDes"
$brand.Range("F1").Value = "Brand Name"
$brand.Range("G1").Value = "Address Line"
$brand.Range("H1").Value = "City"
$brand.Range("I1").Value = "State"
$brand.Range("J1").Value = "Postal"
$brand.Range("K1").Value = "Phone"
$brand.Range("L1").Value = "Reciepents"
$brand.Range("M1").Value = "Sender"
$brand.Range("N1").Value = "Email Subject"
$brand.Range("O1").Value = "Email Banner"
$brand.Range("P1").Value = "Email Body"
$brand.Range("Q1").Value = "Email Signature"
$brand.Range("R1").Value = "1Star Messaging"
$brand.Range("S1").Value = "3Star Messaging"
$brand.Range("T1").Value = "4Star Messaging"
$brand.Range("U1").Value = "Rating"
$brand.Range("V1").Value = "Name"
$brand.Range("W1").Value = "Comments"
$brand.Range("X1").Value = "ReSchedule CampName"

$brand.Range("A1:X1").Font.Bold = $true
$brand.Range("A1:X1").WrapText = $true

# Data row 2
$brand.Range("A2").Value = "Email"
$brand.Range("B2").Value = "Brand"
$brand.Range("C2").Value = "English"
$brand.Range("D2").Value = "Campaign Test"
$brand.Range("E2").Value = "Test "
$brand.Range("F2").Value = "Auto Test"
$brand.Range("G2").Value = "Test"
$brand.Range("H2").Value = "Trivandrum"
$brand.Range("I2").Value = "Kerala"
$brand.Range("J2").Value = "'658881"
$brand.Range("K2").Value = "'9098674532"
$brand.Range("M2").Value = "Avinash"
$brand.Range("N2").Value = "Campaign FeedBack"
$brand.Range("O2").Value = "Please provide feedback"
$brand.Range("P2").Value = "Please feel free to share your experience"
$brand.Range("Q2").Value = "Thanks and Regards" + [char]10 + "Avinash"
$brand.Range("R2").Value = "Sorry for your experience and we'll make sure you will have better experience next time"
$brand.Range("S2").Value = "Thanks you for feedback"
$brand.Range("T2").Value = "Thanks for your feedback! Have a greatday"
$brand.Range("U2").Value = "2"
$brand.Range("V2").Value = "Avinash"
$brand.Range("W2").Value = "I had a poor experience"

$brand.Rows.Item(2).RowHeight = 45

$brand.Range("J2:K2").WrapText = $true

# Hyperlink for L2 (Reciepents) -> mailto, with wrap text
$brand.Hyperlinks.Add($brand.Range("L2"), "mailto:amahadev@dacgroup.com", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "amahadev@dacgroup.com")
$brand.Range("L2").WrapText = $true

$brand.Range("A2:I2").WrapText = $true
$brand.Range("M2:T2").WrapText = $true

$brand.PageSetup.PaperSize = 9
$brand.PageSetup.Orientation = 1

$brand.Range("A1").Select()

# ---------------------------------------------------------------------
# 3. GRLocationCampaign sheet content
# ---------------------------------------------------------------------
$grLoc.Range("A1").Value = "CamType"
$grLoc.Range("B1").Value = "CamOption"
$grLoc.Range("C1").Value = "CamLang"
$grLoc.Range("D1").Value = "CamName"
$grLoc.Range("E1").Value = "CamDes"
$grLoc.Range("F1").Value = "Location"
$grLoc.Range("G1").Value = "1Star Messaging"
$grLoc.Range("H1").Value = "3Star Messaging"
$grLoc.Range("I1").Value = "4Star Messaging"

$grLoc.Range("A2").Value = "General Review Link"
$grLoc.Range("B2").Value = "Location"
$grLoc.Range("C2").Value = "English"
$grLoc.Range("D2").Value = "General R Location"
$grLoc.Range("E2").Value = "Test "
$grLoc.Range("F2").Value = "NTBACF02"
$grLoc.Range("G2").Value = "Sorry for your experience and we'll make sure you will have better experience next time"
$grLoc.Range("H2").Value = "Thanks you for feedback"
$grLoc.Range("I2").Value = "Thanks for your feedback! Have a greatday"

$grLoc.Range("B2").Select()
$grLoc.Activate()

# ---------------------------------------------------------------------
# 4. Update EmailLocationCampaign sheet (S2 value + selection)
# ---------------------------------------------------------------------
$emailLoc.Range("S2").Value = "Campaign Test06/15/2020 8:33:17 PM"

$emailLoc.Activate()
$emailLoc.Range("P1:S2").Select()

# ---------------------------------------------------------------------
# 5. Workbook-level view: active tab = GRLocationCampaign (index 3)
# ---------------------------------------------------------------------
$grLoc.Activate()

Write-Host "Edit complete"
